$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product")

# Update the ID values in column A (rows 3-6) to be zero-indexed
$ws.Range("A3").Value = 0
$ws.Range("A4").Value = 1
$ws.Range("A5").Value = 2
$ws.Range("A6").Value = 3

# Update the selected cell / active cell to A6
$ws.Activate()
$ws.Range("A6").Select()
